# Scheduled-runner style refresh of market-board derived columns
# (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# across the per-job Leve tables. Values below come straight from the
# upstream data refresh; only cells that actually changed are touched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 286
$ws.Range("I9").Value = 309.66666
$ws.Range("J9").Value = 268.25
$ws.Range("K9").Value = 309.66666
$ws.Range("L9").Value = 268.25
$ws.Range("M9").Value = -140.66666
$ws.Range("N9").Value = -606.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 10003.5
$ws.Range("I16").Value = 10003.5
$ws.Range("K16").Value = 10003.5
$ws.Range("M16").Value = -9773.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 1107.4706
$ws.Range("I55").Value = 546
$ws.Range("J55").Value = 2136.8333
$ws.Range("K55").Value = 546
$ws.Range("L55").Value = 2136.8333
$ws.Range("M55").Value = -332
$ws.Range("N55").Value = -2564.8333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 79998
$ws.Range("J87").Value = 79998
$ws.Range("L87").Value = 79998
$ws.Range("N87").Value = -82494

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 79998
$ws.Range("J90").Value = 79998
$ws.Range("L90").Value = 239994
$ws.Range("N90").Value = -252474

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 39768.92
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2866.5
$ws.Range("I113").Value = 2425
$ws.Range("J113").Value = 3749.5
$ws.Range("K113").Value = 2425
$ws.Range("L113").Value = 3749.5
$ws.Range("M113").Value = 829
$ws.Range("N113").Value = -10257.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4151.016
$ws.Range("J138").Value = 6265.087
$ws.Range("L138").Value = 18795.261
$ws.Range("N138").Value = -29075.261

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1337.1333
$ws.Range("I2").Value = 1104.2727
$ws.Range("J2").Value = 1977.5
$ws.Range("K2").Value = 1104.2727
$ws.Range("L2").Value = 1977.5
$ws.Range("M2").Value = -991.2727
$ws.Range("N2").Value = -2203.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1238.8889
$ws.Range("I45").Value = 1238.8889
$ws.Range("K45").Value = 1238.8889
$ws.Range("M45").Value = -861.8888999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 35000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H65").Value = 35000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1337.1333
$ws.Range("I116").Value = 1104.2727
$ws.Range("J116").Value = 1977.5
$ws.Range("K116").Value = 1104.2727
$ws.Range("L116").Value = 1977.5
$ws.Range("M116").Value = 1189.7273
$ws.Range("N116").Value = -6565.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2675.1667
$ws.Range("I122").Value = 2736.5454
$ws.Range("K122").Value = 8209.636200000001
$ws.Range("M122").Value = -5759.636200000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1337.1333
$ws.Range("I3").Value = 1104.2727
$ws.Range("J3").Value = 1977.5
$ws.Range("K3").Value = 1104.2727
$ws.Range("L3").Value = 1977.5
$ws.Range("M3").Value = -990.2727
$ws.Range("N3").Value = -2205.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").ClearContents()
$ws.Range("N44").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 102.25
$ws.Range("I7").Value = 51.142857
$ws.Range("J7").Value = 221.5
$ws.Range("K7").Value = 51.142857
$ws.Range("L7").Value = 221.5
$ws.Range("M7").Value = 61.857143
$ws.Range("N7").Value = -447.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 500
$ws.Range("I45").Value = 500
$ws.Range("K45").Value = 500
$ws.Range("M45").Value = 93

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1980.5938
$ws.Range("I58").Value = 1144.9259
$ws.Range("J58").Value = 6493.2
$ws.Range("K58").Value = 1144.9259
$ws.Range("L58").Value = 6493.2
$ws.Range("M58").Value = -941.9259
$ws.Range("N58").Value = -6899.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1980.5938
$ws.Range("I136").Value = 1144.9259
$ws.Range("J136").Value = 6493.2
$ws.Range("K136").Value = 3434.7777
$ws.Range("L136").Value = 19479.6
$ws.Range("M136").Value = -884.7776999999996
$ws.Range("N136").Value = -24579.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 2254
$ws.Range("I17").Value = 2254
$ws.Range("K17").Value = 2254
$ws.Range("M17").Value = -2086

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 14579.111
$ws.Range("I80").Value = 4651
$ws.Range("J80").Value = 34435.332
$ws.Range("K80").Value = 4651
$ws.Range("L80").Value = 34435.332
$ws.Range("M80").Value = -3653
$ws.Range("N80").Value = -36431.332

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 14579.111
$ws.Range("I83").Value = 4651
$ws.Range("J83").Value = 34435.332
$ws.Range("K83").Value = 23255
$ws.Range("L83").Value = 172176.66
$ws.Range("M83").Value = -18263
$ws.Range("N83").Value = -182160.66

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 7029.8667
$ws.Range("I92").Value = 25000
$ws.Range("K92").Value = 25000
$ws.Range("M92").Value = -23128

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 63426.293
$ws.Range("J122").Value = 130921
$ws.Range("L122").Value = 392763
$ws.Range("N122").Value = -397663

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2915.6667
$ws.Range("I22").Value = 1720.3667
$ws.Range("K22").Value = 1720.3667
$ws.Range("M22").Value = -1425.3667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2915.6667
$ws.Range("I27").Value = 1720.3667
$ws.Range("K27").Value = 1720.3667
$ws.Range("M27").Value = -1613.3667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 50000
$ws.Range("I63").Value = 50000
$ws.Range("K63").Value = 50000
$ws.Range("M63").Value = -49251

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H66").Value = 50000
$ws.Range("I66").Value = 50000
$ws.Range("K66").Value = 150000
$ws.Range("M66").Value = -146256

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 60000
$ws.Range("J133").Value = 60000
$ws.Range("L133").Value = 60000
$ws.Range("N133").Value = -65060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H88").Value = 57500
$ws.Range("J88").Value = 60000
$ws.Range("L88").Value = 60000
$ws.Range("N88").Value = -60812

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H91").Value = 57500
$ws.Range("J91").Value = 60000
$ws.Range("L91").Value = 60000
$ws.Range("N91").Value = -62808

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2125.25
$ws.Range("I107").Value = 2157.4285
$ws.Range("K107").Value = 6472.2855
$ws.Range("M107").Value = -4552.2855
